$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44304, 1, 8, 86.16975441619991),
    @(44305, 2, 9, 96.9409737182249),
    @(44306, 2, 11, 118.4834123222749),
    @(44307, 0, 10, 107.7121930202499)
)

$styleSource = $ws.Cells.Item(229, 1)
$startRow = 230

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]

    # Match the date-column cell style (border/font/alignment/number format)
    # used by the rest of column A, like Excel would when extending a series.
    $styleSource.Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
